$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Arquivo alterado"
$ws.Range("A2").Value = "ok"

$ws.Range("B2").Select()
